$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Control 25
$ws.Range("D2").Value = 0.2492544832245027
$ws.Range("E2").Value = 0.2492544832245027

# Row 3 - Control 44
$ws.Range("D3").Value = [double]"2.712257442123737E-21"
$ws.Range("E3").Value = [double]"2.712257442123737E-21"

# Row 4 - Control 40
$ws.Range("D4").Value = 0.000135036822210734
$ws.Range("E4").Value = 0.000135036822210734

# Row 5 - Control 41
$ws.Range("D5").Value = 0.16236598794312
$ws.Range("E5").Value = 0.16236598794312

# Row 7 - MDD 38
$ws.Range("D7").Value = 0.9980256021924404
$ws.Range("E7").Value = 0.001974397807559591

# Row 8 - MDD 9
$ws.Range("D8").Value = 0.9999999995418654
$ws.Range("E8").Value = [double]"4.581346413345955E-10"

# Row 9 - MDD 49
$ws.Range("D9").Value = 0.1162655217639005
$ws.Range("E9").Value = 0.8837344782360995

# Row 10 - MDD 26
$ws.Range("D10").Value = 0.9999999999999991
$ws.Range("E10").Value = [double]"8.881784197001252E-16"

# Row 11 - MDD 34
$ws.Range("D11").Value = 0.0230164776105816
$ws.Range("E11").Value = 0.9769835223894184
$ws.Range("F11").Value = 4.578420639038086
